$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.739.48"
$ws.Range("E2").Value = "  -0.31%  "

$ws.Range("D3").Value = "2.677.63"
$ws.Range("E3").Value = "  -1.11%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.95"
$ws.Range("E5").Value = "  -1.38%  "

$ws.Range("E6").Value = "  -0.96%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  +4.96%  "

$ws.Range("E9").Value = "  +4.21%  "

$ws.Range("E10").Value = "  -1.88%  "

$ws.Range("E12").Value = "  -0.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.28"
$ws.Range("E13").Value = "  -4.10%  "

$ws.Range("E14").Value = "  -3.20%  "

$ws.Range("D15").Value = "3.158.56"
$ws.Range("E15").Value = "  -1.08%  "

$ws.Range("D16").Value = "65.601.98"

$ws.Range("D17").Value = "2.678.60"
$ws.Range("E17").Value = "  -1.19%  "

$ws.Range("E18").Value = "  +1.43%  "

$ws.Range("E19").Value = "  -2.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.58"
$ws.Range("E20").Value = "  -1.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.94"
$ws.Range("E21").Value = "  -2.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.89"
$ws.Range("E23").Value = "  -1.49%  "

$ws.Range("E24").Value = "  +5.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.66"
$ws.Range("E25").Value = "  -2.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.66"
$ws.Range("E26").Value = "  -0.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.167"
$ws.Range("E27").Value = "  -2.17%  "

$ws.Range("E28").Value = "  -6.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.05"
$ws.Range("E29").Value = "  -4.36%  "

$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "533.90"
$ws.Range("E31").Value = "  -2.14%  "

$ws.Range("E32").Value = "  -3.14%  "

$ws.Range("E33").Value = "  -2.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.55"
$ws.Range("E34").Value = "  +1.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.50"
$ws.Range("E35").Value = "  -4.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.424"
$ws.Range("E36").Value = "  -2.60%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.60"
$ws.Range("E37").Value = "  -1.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "159.35"
$ws.Range("E38").Value = "  -2.52%  "

$ws.Range("E39").Value = "  -0.05%  "

$ws.Range("E40").Value = "  -2.66%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "163.87"
$ws.Range("E42").Value = "  -5.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.14"
$ws.Range("E43").Value = "  -1.71%  "

$ws.Range("E44").Value = "  +1.84%  "

$ws.Range("E45").Value = "  -1.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.84"
$ws.Range("E46").Value = "  -3.27%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0258"
$ws.Range("E47").Value = "  -3.33%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.640"
$ws.Range("E48").Value = "  -2.50%  "

$ws.Range("E49").Value = "  +14.88%  "

$ws.Range("E50").Value = "  -4.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0998"
$ws.Range("E51").Value = "  +0.59%  "
